# Generate Report for Handoff
# Adds two new "ready for handoff" rows (57c5e2f7-... and 5c1c155d-...) to the
# Overview sheet and to each language sheet (zh-cn, de-de), mirroring the
# shape of the existing rows already on each sheet.

$wb = $excel.ActiveWorkbook

$guid1 = "57c5e2f7-6f2c-4760-af75-15273759b289"
$guid2 = "5c1c155d-8154-4fe3-81aa-4ccca6c912af"

$hash1 = "427ba13591edbe2559cba2f27b808780a1142311"
$hash2 = "321f26decdcc42e6762a4724773ddf6576573b6c"

$statusReady = "Ready for handoff"
$noHandback  = "0001-01-01 00:00:00"
$reason      = "Include"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(4, 1), "https://github.com/OpenLocalizationTest/oltest/blob/" + $hash1 + "/e2e/" + $guid1 + ".md", $null, $null, $guid1 + ".md")
$wsOverview.Cells.Item(4, 2).Value = $statusReady
$wsOverview.Cells.Item(4, 3).Value = $statusReady
$wsOverview.Cells.Item(4, 4).Value = "2016-29-20 02:29:23"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(5, 1), "https://github.com/OpenLocalizationTest/oltest/blob/" + $hash2 + "/e2e/" + $guid2 + ".md", $null, $null, $guid2 + ".md")
$wsOverview.Cells.Item(5, 2).Value = $statusReady
$wsOverview.Cells.Item(5, 3).Value = $statusReady
$wsOverview.Cells.Item(5, 4).Value = "2016-29-20 02:29:23"

# ---------------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | Latest Target File | Latest Handback File
# | Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4, 1), "https://github.com/OpenLocalizationTest/oltest/blob/" + $hash1 + "/e2e/" + $guid1 + ".md", $null, $null, $guid1 + ".md")
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4, 2), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/" + $hash1 + "/e2e/" + $guid1 + ".md", $null, $null, ".md")
$wsZh.Cells.Item(4, 3).Value = $statusReady
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(4, 4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $hash1 + "/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $guid1 + "." + $hash1 + ".zh-cn.xlf", $null, $null, $guid1 + "." + $hash1 + ".zh-cn.xlf")
$wsZh.Cells.Item(4, 5).Value = "2016-03-20 02:29:20"
$wsZh.Cells.Item(4, 8).Value = $noHandback
$wsZh.Cells.Item(4, 9).Value = $reason

$wsZh.Hyperlinks.Add($wsZh.Cells.Item(5, 1), "https://github.com/OpenLocalizationTest/oltest/blob/" + $hash2 + "/e2e/" + $guid2 + ".md", $null, $null, $guid2 + ".md")
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(5, 2), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/" + $hash2 + "/e2e/" + $guid2 + ".md", $null, $null, ".md")
$wsZh.Cells.Item(5, 3).Value = $statusReady
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(5, 4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $hash2 + "/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $guid2 + "." + $hash2 + ".zh-cn.xlf", $null, $null, $guid2 + "." + $hash2 + ".zh-cn.xlf")
$wsZh.Cells.Item(5, 5).Value = "2016-03-20 02:29:20"
$wsZh.Cells.Item(5, 8).Value = $noHandback
$wsZh.Cells.Item(5, 9).Value = $reason

# ---------------------------------------------------------------------------
# de-de sheet: same shape as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4, 1), "https://github.com/OpenLocalizationTest/oltest/blob/" + $hash1 + "/e2e/" + $guid1 + ".md", $null, $null, $guid1 + ".md")
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4, 2), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/" + $hash1 + "/e2e/" + $guid1 + ".md", $null, $null, ".md")
$wsDe.Cells.Item(4, 3).Value = $statusReady
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(4, 4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $hash1 + "/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $guid1 + "." + $hash1 + ".de-de.xlf", $null, $null, $guid1 + "." + $hash1 + ".de-de.xlf")
$wsDe.Cells.Item(4, 5).Value = "2016-03-20 02:29:23"
$wsDe.Cells.Item(4, 8).Value = $noHandback
$wsDe.Cells.Item(4, 9).Value = $reason

$wsDe.Hyperlinks.Add($wsDe.Cells.Item(5, 1), "https://github.com/OpenLocalizationTest/oltest/blob/" + $hash2 + "/e2e/" + $guid2 + ".md", $null, $null, $guid2 + ".md")
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(5, 2), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/" + $hash2 + "/e2e/" + $guid2 + ".md", $null, $null, ".md")
$wsDe.Cells.Item(5, 3).Value = $statusReady
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(5, 4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $hash2 + "/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $guid2 + "." + $hash2 + ".de-de.xlf", $null, $null, $guid2 + "." + $hash2 + ".de-de.xlf")
$wsDe.Cells.Item(5, 5).Value = "2016-03-20 02:29:23"
$wsDe.Cells.Item(5, 8).Value = $noHandback
$wsDe.Cells.Item(5, 9).Value = $reason
